$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "27.138.64"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "1.566.16"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  +0.18%  "
Set-TextValue "D5" "208.03"
$ws.Range("E5").Value = "  +1.00%  "
Set-TextValue "D6" "0.492"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("E7").Value = "  +0.10%  "
Set-TextValue "D8" "22.23"
$ws.Range("E8").Value = "  +4.31%  "
$ws.Range("E9").Value = "  +0.99%  "
Set-TextValue "D10" "0.0588"
$ws.Range("E10").Value = "  +1.22%  "
Set-TextValue "D11" "0.0862"
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").Value = "1.789.83"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").Value = "1.566.28"
$ws.Range("E13").Value = "  +1.61%  "
Set-TextValue "D14" "3.76"
$ws.Range("E14").Value = "  +2.30%  "
Set-TextValue "D15" "0.522"
$ws.Range("E15").Value = "  +2.46%  "
$ws.Range("D16").Value = "27.125.77"
$ws.Range("E16").Value = "  +1.23%  "
Set-TextValue "D17" "62.09"
$ws.Range("E17").Value = "  +1.58%  "
Set-TextValue "D18" "219.94"
$ws.Range("E18").Value = "  +2.86%  "
$ws.Range("E19").Value = "  +2.37%  "
Set-TextValue "D20" "7.35"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("E21").Value = "  +0.20%  "
Set-TextValue "D22" "4.08"
$ws.Range("E22").Value = "  +1.92%  "
Set-TextValue "D23" "9.30"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("E24").Value = "  +1.51%  "
Set-TextValue "D25" "154.41"
$ws.Range("E25").Value = "  +1.80%  "
Set-TextValue "D26" "6.64"
$ws.Range("E26").Value = "  +0.72%  "
Set-TextValue "D27" "15.04"
$ws.Range("E27").Value = "  +1.67%  "
Set-TextValue "D28" "1.00"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("E29").Value = "  +1.72%  "
Set-TextValue "D30" "0.0470"
$ws.Range("E30").Value = "  +2.49%  "
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").Value = "1.445.94"
$ws.Range("E33").Value = "  +5.85%  "
$ws.Range("E34").Value = "  +5.27%  "
Set-TextValue "D35" "1.58"
$ws.Range("E35").Value = "  +4.57%  "
Set-TextValue "D36" "0.969"
$ws.Range("E36").Value = "  +0.83%  "
Set-TextValue "D37" "2.29"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("E39").Value = "  +0.89%  "
Set-TextValue "D40" "0.816"
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("E43").Value = "  +3.72%  "
Set-TextValue "D44" "0.989"
$ws.Range("E44").Value = "  -0.23%  "
Set-TextValue "D45" "64.61"
$ws.Range("E45").Value = "  +2.59%  "
Set-TextValue "D46" "1.78"
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("D47").Value = "1.703.83"
$ws.Range("E47").Value = "  +1.84%  "
Set-TextValue "D48" "87.04"
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D49" "0.0524"
$ws.Range("E49").Value = "  +3.00%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0102"
$ws.Range("E50").Value = "  +4.60%  "
Set-TextValue "D51" "0.0967"
$ws.Range("E51").Value = "  +2.52%  "
